$wb = $excel.ActiveWorkbook

# --- Sheet "Info": update Objetivo/Tiempo result row ---
$wsInfo = $wb.Worksheets.Item("Info")
$wsInfo.Range("A2").Value = 640108574274.0112
$wsInfo.Range("B2").Value = 2.139999866485596

# --- Sheet "Activados": Proceso column becomes 1, Tiempo extended 0..360 step 20 (rows 2-20) ---
$wsAct = $wb.Worksheets.Item("Activados")
for ($i = 2; $i -le 20; $i++) {
    $wsAct.Cells.Item($i, 1).Value = 1
    $wsAct.Cells.Item($i, 2).Value = ($i - 2) * 20
}

# --- Sheet "Operando": Proceso column changes from 4 to 1 for every existing row (2-366) ---
$wsOp = $wb.Worksheets.Item("Operando")
for ($i = 2; $i -le 366; $i++) {
    $wsOp.Cells.Item($i, 1).Value = 1
}

# --- Sheet "Contaminantes": update Z and Concentracion columns for rows 2-6 ---
$wsCont = $wb.Worksheets.Item("Contaminantes")
$wsCont.Cells.Item(2, 2).Value = 449208244800.0004
$wsCont.Cells.Item(2, 3).Value = 16.66000000000001
$wsCont.Cells.Item(3, 2).Value = 13481640000.00001
$wsCont.Cells.Item(3, 3).Value = 0.5000000000000004
$wsCont.Cells.Item(4, 2).Value = 87091394399.99998
$wsCont.Cells.Item(4, 3).Value = 3.23
$wsCont.Cells.Item(5, 2).Value = 307074.010608
$wsCont.Cells.Item(5, 3).Value = 0.0000113886
$wsCont.Cells.Item(6, 2).Value = 90326988000.00008
$wsCont.Cells.Item(6, 3).Value = 3.350000000000003
